# Fix test-case name error: the daily report was missing the row for the
# 17-10-31 ~ 17-11-01 period. Append that row (row 37) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset to the default "Normal" style first so the new cells don't inherit
# column-level formatting (e.g. column L's wrap-text style) - the source
# report writes these rows as plain, unstyled text cells.
$ws.Range("A37:N37").Style = "Normal"

# Plain text value for A37 (not number-like, so Excel's input parser
# leaves it as literal text without any extra styling).
$ws.Range("A37").Value = "17-10-31 ~ 17-11-01"

# Number-LIKE text values ("80", "5", "0", "6.25%", "100.00%", "0.00%", ...)
# must stay as literal text, matching the rest of the sheet. A direct
# .Value assignment would be auto-converted to a real number by Excel's
# input parser, so build each as a text-returning formula first and then
# convert that formula to a static value via copy / paste-special-values,
# which preserves the literal text without stamping a quote-prefix /
# text-number-format style onto the cell.
$ws.Range("B37").Formula = '="80"'
$ws.Range("C37").Formula = '="5"'
$ws.Range("D37").Formula = '="6.25%"'
$ws.Range("E37").Formula = '="0"'
$ws.Range("F37").Formula = '="5"'
$ws.Range("G37").Formula = '="100.00%"'
$ws.Range("H37").Formula = '="5"'
$ws.Range("I37").Formula = '="0"'
$ws.Range("J37").Formula = '="0"'
$ws.Range("K37").Formula = '="0.00%"'

$ws.Range("B37:K37").Copy()
$ws.Range("B37").PasteSpecial(-4163)

# Remaining plain text values, set after B:K so new shared-string entries
# land in the same order the source workbook registers them in (A37's date
# label, then D37's "6.25%", then M37's device list).
$ws.Range("L37").Value = "[]"
$ws.Range("M37").Value = "['7ff9010202000024', '7ff9010202000025', '7ff9010202000026', '7ff9010202000027', '7ff9010202000028']"
$ws.Range("N37").Value = "[]"

$excel.CutCopyMode = $false
